# Gantt chart update for 6/21/2020
# - Duplicate sheet "200515" into a new sheet "200621" placed right after it
# - Tweak a couple of dates/labels that changed since the last snapshot
# - Append three new tasks (rows 18-20) to the new sheet
# - New sheet becomes the active / selected tab; old sheet keeps a plain
#   selection over its data range and loses the "tabSelected" flag

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

$old = $wb.Worksheets.Item("200515")
$old.Copy($null, $old)

$new = $wb.Worksheets.Item($old.Index + 1)
$new.Name = "200621"

# Seed the new task-name strings in the same order they were typed while
# building the updated chart, so they land in the shared string table the
# same way they did originally.
$new.Range("A20").Value = "Continue to work on project for publication"
$new.Range("A16").Value = "Find experimental data values for outputs"
$new.Range("A19").Value = "Build prediction tool"
$new.Range("A17").Value = "Descriptor expansion/selection through SISSO and LASSO"

# --- Row 13: "Uncertainty prediction: all models" end date pushed out ---
$new.Range("D13").Value = [DateTime]"2020-06-15"

# --- Row 16: used to be "Expected improvement..." - now a different task ---
$new.Range("B16").Value = "Lit Review"
$new.Range("C16").Value = [DateTime]"2020-06-10"
$new.Range("D16").Value = [DateTime]"2020-07-01"

# --- Row 17: used to be "Have model completed" - now a different task ---
$new.Range("B17").Value = "Total data"
$new.Range("C17").Value = [DateTime]"2020-05-25"
$new.Range("D17").Value = [DateTime]"2020-06-15"

$new.Range("A5").Copy()
$new.Range("A17").PasteSpecial($xlPasteFormats)
$new.Range("C2").Copy()
$new.Range("B17").PasteSpecial($xlPasteFormats)
$new.Range("D7").Copy()
$new.Range("D17").PasteSpecial($xlPasteFormats)

# --- Row 18: "Have model completed" moved down here with new dates ---
$new.Range("A18").Value = "Have model completed"
$new.Range("B18").Value = "Overall"
$new.Range("C18").Value = [DateTime]"2020-06-10"
$new.Range("D18").Value = [DateTime]"2020-06-17"

$new.Range("C2").Copy()
$new.Range("C18:D18").PasteSpecial($xlPasteFormats)

# --- Row 19: brand new task ---
$new.Range("B19").Value = "Overall"
$new.Range("C19").Value = [DateTime]"2020-06-15"
$new.Range("D19").Value = [DateTime]"2020-06-24"

$new.Range("A5").Copy()
$new.Range("A19").PasteSpecial($xlPasteFormats)
$new.Range("C2").Copy()
$new.Range("B19:D19").PasteSpecial($xlPasteFormats)

# --- Row 20: brand new task ---
$new.Range("B20").Value = "Bonus"
$new.Range("C20").Value = [DateTime]"2020-06-24"
$new.Range("D20").Value = [DateTime]"2020-09-01"

$new.Range("A5").Copy()
$new.Range("A20").PasteSpecial($xlPasteFormats)
$new.Range("D7").Copy()
$new.Range("C20:D20").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- sheet views / selection ---
$old.Range("A1:D17").Select()

$new.Range("D21").Select()
$new.Activate()

$wb.Save()
